$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Select()
$ws.Range("C4").Value = "Einführung – Python-Implementierung von Graphen-Suchroblemen"
